$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# foxBMS-slave BOM update to version 2.02
# ---------------------------------------------------------------------------

# 1) The "1uF / C501, C503" capacitor line (row 9) is no longer part of the
#    BOM - remove the whole row and shift everything below it up.
$ws.Rows(9).Delete()

# 2) R302,R304,R305,R310,R402,R404,R405,R410 value changes from the numeric
#    "0" to the textual "0R".
$ws.Range("C6").Value = "0R"

# 3) The transformer part (row 12 after the shift above) changes from
#    HX1188NL to HM2102NL; it no longer has a Digikey order number.
$ws.Range("C12").Value = "HM2102NL"
$ws.Range("D12").Value = "HM2102NL"
$ws.Range("E12").Value = "HM2102NL"
$ws.Range("I12").Clear()

# 4) The isoSPI transceiver (row 13 after the shift above) is now referenced
#    by its generic device name "LTC6820" instead of the exact order code.
$ws.Range("C13").Value = "LTC6820"
$ws.Range("D13").Value = "LTC6820"

# ---------------------------------------------------------------------------
# View / selection bookkeeping that Excel recorded when the file was saved.
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Left = 360
$wb.Windows.Item(1).Top = 360

$ws.Range("F20").Select() | Out-Null
